# Applies the commit "added required experimental boolean element to valuesets"
#
# The "Metadata" worksheet holds a Property/Value table. The "Experimental"
# row (row 7) previously had an empty Value cell; this adds the required
# "true" value (stored as plain text, same as the other Property/Value
# entries on this sheet - not as a native Excel boolean). The sheet's Date
# row (row 8) is also refreshed to the timestamp recorded at the time this
# boolean was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$target = $ws.Range("B7")
$formatSource = $ws.Range("B6")

# Write the literal text "true" into B7. A bare Value = "true" gets
# auto-coerced by Excel into a native boolean TRUE, so force text by typing
# it with a leading apostrophe (COM equivalent of typing '\true\ in the UI).
$target.Value = "'true"

# The apostrophe entry leaves a quote-prefix flag on the cell's style (a
# new cellXf), which the source workbook does not have - B7 keeps the same
# style as the rest of the column (s="2"). Re-apply that original
# formatting by copying it from a neighboring cell that still has it.
$formatSource.Copy() | Out-Null
$target.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Date value refreshed alongside the new Experimental flag
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
